$wb = $excel.ActiveWorkbook

# --- Transactions sheet: record its existing selection before we move away
#     from it (it will no longer be the active/shown tab afterwards) ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("I13:J13").Select() | Out-Null

# --- Repayment Schedule sheet: a new (blank) column is inserted before
#     column N, pushing the existing "Late" / "Outstanding" / "Over Due"
#     columns one position to the right ---
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Columns("N:N").Insert() | Out-Null

# Make "Repayment Schedule" the active sheet/tab with its new selection
$wsRepay.Select()
$wsRepay.Range("I18").Select() | Out-Null
